# Updates the "Periodo Mora" labels (column E) and "Valor Mora" amounts
# (column F) for rows 16-22 on Hoja1 so the periods now run newest -> oldest
# (2404, 2403, 2402, 2401, 2312, 2311, 2310) instead of oldest -> newest,
# with the partial payment amount (153334) following period 2404 and the
# other periods carrying the full amount (200000).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periods = @("2404", "2403", "2402", "2401", "2312", "2311", "2310")
$amounts = @(153334, 200000, 200000, 200000, 200000, 200000, 200000)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $amounts[$i]
}
